$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B37: end time corrected ---
$ws.Range("B37").Value = "12/01/2026 11:54"

# --- Append new occurrence rows 40-44 ---
# Copy formatting from the last existing data row (39) down to each new row
# first, then fill in the values, so styles (s=3 / s=4 for the ALA column)
# match the rest of the table.

$newRows = @(
    @{ Row = 40; A = "12/01/2026 13:06"; B = "12/01/2026 17:38"; C = "4ª"; D = "P"; E = "2026-43493193-0"; F = "P01001 - VISTORIA DE FISCALIZACAO"; G = "RUA AURORA, Nº 17 - SAO BENEDITO - PASSOS"; H = "APV06243" },
    @{ Row = 41; A = "12/01/2026 14:13"; B = "12/01/2026 14:43"; C = "4ª"; D = "V"; E = "2026-43494141-1"; F = "V02201 - VITIMA DE CHOQUE DE BICICLETA"; G = "AVENIDA ESTACAO, Nº 210 - PARQUE DA ESTACAO - PASSOS"; H = "UR04360" },
    @{ Row = 42; A = "12/01/2026 14:23"; B = "12/01/2026 15:32"; C = "4ª"; D = "V"; E = "2026-43494289-2"; F = "V01008 - VITIMA COM DOR ABDOMINAL"; G = "RUA MISSOES, Nº 565 - JARDIM PLANALTO - PASSOS"; H = "UR04360" },
    @{ Row = 43; A = "12/01/2026 17:46"; B = "12/01/2026 18:51"; C = "4ª"; D = "V"; E = "2026-43498275-0"; F = "V01003 - VITIMA COM CRISE CONVULSIVA"; G = "RUA SAGU, Nº 16 - RESIDENCIAL PORTAL DAS PALMEIRAS - PASSOS"; H = "UR04360" },
    @{ Row = 44; A = "12/01/2026 19:39"; B = "12/01/2026 19:42"; C = "4ª"; D = "V"; E = "2026-43500154-1"; F = "V01004 - VITIMA COM CRISE DIABETICA / HIPOGLICEMIA"; G = "RUA ELZO CALIXTO MATTAR, Nº 701 - SANTA RITA - PASSOS"; H = "UR04360" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $srcRow = $rowNum - 1

    $ws.Range("A${srcRow}:H${srcRow}").Copy()
    $ws.Range("A${rowNum}:H${rowNum}").PasteSpecial(-4122)

    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("E$rowNum").Value = $r.E
    $ws.Range("F$rowNum").Value = $r.F
    $ws.Range("G$rowNum").Value = $r.G
    $ws.Range("H$rowNum").Value = $r.H
}

$excel.CutCopyMode = $false
